# Re-run / refresh of the logistic-regression results workbook.
# Updates:
#   - Stats!C9            (mean of Direct_tilt) recomputed
#   - VIF!C4,C6:C14        VIF values recomputed (tiny float drift)
#   - DTR-Features!A2:C13  Decision-tree feature-importance table re-sorted
#                          with refreshed coefficients

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Stats sheet: single mean value nudged
# ---------------------------------------------------------------------
$wsStats = $wb.Worksheets.Item("Stats")
$wsStats.Range("C9").Value = 21.322

# ---------------------------------------------------------------------
# VIF sheet: refreshed variance-inflation-factor values
# ---------------------------------------------------------------------
$wsVif = $wb.Worksheets.Item("VIF")
$wsVif.Range("C4").Value  = 19.55932651283902
$wsVif.Range("C6").Value  = 40.25515896510461
$wsVif.Range("C7").Value  = 2.739011001355327
$wsVif.Range("C8").Value  = 3.819412339414832
$wsVif.Range("C9").Value  = 7.346293972889505
$wsVif.Range("C10").Value = 15.2433120509357
$wsVif.Range("C11").Value = 17.04711841146865
$wsVif.Range("C12").Value = 2.338287830692729
$wsVif.Range("C13").Value = 6.960206749404251
$wsVif.Range("C14").Value = 2.169385148926902

# ---------------------------------------------------------------------
# DTR-Features sheet: feature importances re-sorted after model re-fit.
# Column A = original feature index, B = feature name, C = coefficient
# (stored as text, matching the source workbook).
# ---------------------------------------------------------------------
$wsDtr = $wb.Worksheets.Item("DTR-Features")

$rows = @(
    @{ Row = 2;  A = 5;  B = "degree_spondylolisthesis"; C = "0.535" },
    @{ Row = 3;  A = 3;  B = "sacral_slope";              C = "0.151" },
    @{ Row = 4;  A = 4;  B = "pelvic_radius";              C = "0.095" },
    @{ Row = 5;  A = 9;  B = "cervical_tilt";              C = "0.07"  },
    @{ Row = 6;  A = 1;  B = "pelvic_tilt";                C = "0.037" },
    @{ Row = 7;  A = 10; B = "sacrum_angle";               C = "0.031" },
    @{ Row = 8;  A = 6;  B = "pelvic_slope";               C = "0.027" },
    @{ Row = 9;  A = 0;  B = "pelvic_incidence";           C = "0.019" },
    @{ Row = 10; A = 8;  B = "thoracic_slope";             C = "0.019" },
    @{ Row = 11; A = 11; B = "scoliosis_slope";            C = "0.016" },
    @{ Row = 12; A = 2;  B = "lumbar_lordosis_angle";      C = "0.0"   },
    @{ Row = 13; A = 7;  B = "Direct_tilt";                C = "0.0"   }
)

foreach ($r in $rows) {
    $wsDtr.Range("A$($r.Row)").Value = $r.A
    $wsDtr.Range("B$($r.Row)").Value = $r.B

    $existingText = $wsDtr.Range("C$($r.Row)").Text
    if ($existingText -ne $r.C) {
        # Leading apostrophe forces text storage so "0.070" style values
        # keep their exact textual form (matches source: column C is text,
        # not a number) instead of being normalised by Excel's numeric
        # auto-detect. Only touch cells whose text actually changes so
        # untouched rows (e.g. the repeated "0.0" rows) keep their original
        # formatting untouched.
        $wsDtr.Range("C$($r.Row)").Formula = "'" + $r.C
    }
}
